$p = $ppt.ActivePresentation

# 1. Merge the three runs of the second paragraph in the Content Placeholder
#    on slide 9 into a single run (the "and popularity" split is removed).
$s9 = $p.Slides.Item(9)
$shp = $s9.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$para2 = $tr.Paragraphs(2, 1)
# Force a real text change first (the host no-ops a set that is byte-identical
# to the existing concatenated text), then apply the final merged text so the
# three runs collapse into one run using the first run's formatting.
$para2.Text = "__TMP_PLACEHOLDER__"
$para2b = $tr.Paragraphs(2, 1)
$para2b.Text = "We also wanted to analyze which part of town needs more parking spaces based on the availability and popularity of the current parking meters."

# 2. Delete the last (blank) slide from the deck.
$p.Slides.Item($p.Slides.Count).Delete()
